$wb = $excel.ActiveWorkbook

# --- Reorder worksheet tabs: "Konfiguration" moves in front of "Beitraege" ---
# (Gebuehren stays last / unaffected)
$wb.Worksheets.Item("Konfiguration").Move($wb.Worksheets.Item(1))

# Re-resolve the worksheet reference by name now that the tab order changed -
# stale references captured before .Move() keep pointing at the old position.
$konfig = $wb.Worksheets.Item("Konfiguration")

# --- Add a new configuration row: "Rechnungsstartindex" = 1 ---
$konfig.Range("A3").Value = "Rechnungsstartindex"
$konfig.Range("B3").Value = 1

# Widen column A so the new, longer label fits.
$konfig.Columns.Item(1).ColumnWidth = 17.7

# Keep the selection / active cell on the newly added row.
[void]$konfig.Range("A3").Select()

# "Konfiguration" becomes the active (selected) sheet tab.
[void]$konfig.Activate()
